# Bala Added Google sign up test cases
$wb = $excel.ActiveWorkbook
$searchSheet = $wb.ActiveSheet

# Insert the new sheet right after "Search" and name it.
$ws = $wb.Worksheets.Add($null, $searchSheet)
$ws.Name = "GoogleSignupData"

# ---- Header row (row 1) ----
$ws.Cells.Item(1, 1).Value = "FirstName"
$ws.Cells.Item(1, 2).Value = "LastName"
$ws.Cells.Item(1, 3).Value = "email"
$ws.Cells.Item(1, 4).Value = "CreatePassword"
$ws.Cells.Item(1, 5).Value = "ConfirmPassword"
$ws.Cells.Item(1, 6).Value = "BirthMonth"
$ws.Cells.Item(1, 7).Value = "BirthDay"
$ws.Cells.Item(1, 8).Value = "Birthyear"
$ws.Cells.Item(1, 9).Value = "Gender"
$ws.Cells.Item(1, 10).Value = "Moblienumber"
$ws.Cells.Item(1, 11).Value = "currentemail"
$ws.Cells.Item(1, 12).Value = "location"

# ---- Data row (row 2) ----
$ws.Cells.Item(2, 1).Value = "Selenium"
$ws.Cells.Item(2, 2).Value = "Training"
$ws.Cells.Item(2, 3).Value = "selenium.trainingsridhar@gmail.com"
$ws.Cells.Item(2, 4).Value = "123@Pass"
$ws.Cells.Item(2, 5).Value = "123@Pass"
$ws.Cells.Item(2, 6).Value = "May"
$ws.Cells.Item(2, 7).Value = 2
$ws.Cells.Item(2, 8).Value = 2000
$ws.Cells.Item(2, 9).Value = "Male"
$ws.Cells.Item(2, 10).Value = 9876543212
$ws.Cells.Item(2, 11).Value = "****@gmail.com"
$ws.Cells.Item(2, 12).Value = "India"

# ---- Hyperlinks (Excel's "autoformat as you type" turned the "@"
# containing entries into mailto links; this also applies the
# built-in Hyperlink cell style). ----
$ws.Hyperlinks.Add($ws.Cells.Item(2, 3), "mailto:selenium.trainingsridhar@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(2, 4), "mailto:123@Pass") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(2, 11), "mailto:****@gmail.com") | Out-Null

# Leave the new sheet as the active one, with A4 selected.
$ws.Range("A4").Select() | Out-Null
